# trinh | navigate to Category page from Homepage
# Update the "Login" sheet's sample data row (row 2) with a fresh set of
# generated registration data (username/password/name/re_password/phone/birthday).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Plain text values (safe to assign directly - Excel will not try to
# reinterpret them as numbers).
$ws.Range("A2").Value = "user98348@gmail.com"
$ws.Range("B2").Value = "Aty8209@"
$ws.Range("C2").Value = "ziyauxsi"
$ws.Range("D2").Value = "Aty8209@"

# Numeric-looking values (phone number / birthday) must stay text so that
# leading zeros are preserved, exactly like the existing rows. Writing
# them straight into .Value would make Excel coerce them into numbers,
# so instead we build a text value via a formula in a scratch cell and
# paste only the resulting value (which keeps it typed as text) into the
# target cell, then clean up the scratch cell.
$ws.Range("H1").Formula = "=""0975392963"""
$ws.Range("H1").Copy()
$ws.Range("E2").PasteSpecial(-4163)

$ws.Range("H1").Formula = "=""14121988"""
$ws.Range("H1").Copy()
$ws.Range("F2").PasteSpecial(-4163)

$ws.Range("H1").Value = ""
$excel.CutCopyMode = 0
